# Weekly update: new price-report week added for
# Hortaliza, Agrícola del Norte S.A. de Arica - Zapallo italiano.
#
# A new pair of rows (Primera / Segunda calidad) for the latest reporting
# date (2021-10-20, serial 44489) is inserted at the top of the weekly
# data block (row 153), pushing the existing historical rows down by two
# rows. The two oldest rows that fall off the bottom of the table are
# appended again at the end (rows 217/218), matching how this workbook's
# rolling weekly log is maintained.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 153 - this shifts the existing rows
# 153:216 down to 155:218 (including their formatting), which already
# reproduces rows 217/218 (old 215/216) exactly as required.
$ws.Rows.Item(153).Insert()
$ws.Rows.Item(153).Insert()

# Fill the two newly-inserted blank rows with the new week's data - same
# as what used to be in rows 153/154, except for the updated date.
$ws.Range("A153").Value = 1
$ws.Range("B153").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C153").Value = "Arica y Parinacota"
$ws.Range("D153").Value = 44489
$ws.Range("E153").Value = 15
$ws.Range("F153").Value = 100112032
$ws.Range("G153").Value = "Zapallo italiano"
$ws.Range("H153").Value = "Huracán"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 120
$ws.Range("K153").Value = 7000
$ws.Range("L153").Value = 8000
$ws.Range("M153").Value = 7500
$ws.Range("N153").Value = "$/caja 70 unidades"
$ws.Range("O153").Value = "Región de Arica y Parinacota"
$ws.Range("P153").Value = 107
$ws.Range("Q153").Value = 70
$ws.Range("R153").Value = "Hortaliza"

$ws.Range("A154").Value = 1
$ws.Range("B154").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C154").Value = "Arica y Parinacota"
$ws.Range("D154").Value = 44489
$ws.Range("E154").Value = 15
$ws.Range("F154").Value = 100112032
$ws.Range("G154").Value = "Zapallo italiano"
$ws.Range("H154").Value = "Huracán"
$ws.Range("I154").Value = "Segunda"
$ws.Range("J154").Value = 120
$ws.Range("K154").Value = 6000
$ws.Range("L154").Value = 7000
$ws.Range("M154").Value = 6500
$ws.Range("N154").Value = "$/caja 100 unidades"
$ws.Range("O154").Value = "Región de Arica y Parinacota"
$ws.Range("P154").Value = 65
$ws.Range("Q154").Value = 100
$ws.Range("R154").Value = "Hortaliza"
